# Unit Tests for XLSX File Upload (#140)
# Replace the placeholder "888-888-XXXX" style fake phone numbers in the
# "Individuals" sheet with realistic looking phone numbers. The values
# alternate row by row between two phone numbers for each of the two
# phone columns (H = phone_number_1, I = phone_number_2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

$phone1Odd  = "+44 1632 960852"
$phone1Even = "+1-613-555-0182"
$phone2Odd  = "+1-541-754-3010"
$phone2Even = "+36 55 979 922"

for ($row = 3; $row -le 29; $row++) {
    if (($row % 2) -eq 1) {
        $ws.Cells.Item($row, 8).Value = $phone1Odd
        $ws.Cells.Item($row, 9).Value = $phone2Odd
    } else {
        $ws.Cells.Item($row, 8).Value = $phone1Even
        $ws.Cells.Item($row, 9).Value = $phone2Even
    }
}
